# Fix ShortCircuitTest conversion for 3-winding transformers
# - Transformer sheet: correct leakage impedance figures for Reg/tpoletop
#   3-winding (star-point) transformer rows, and rename the shared star bus
#   from "xfstar_tpoletop_B" to "xfbus_tpoletop_B".
# - Bus sheet: the star-bus rename shifts the A/B/C phase rotation for a
#   number of downstream buses; update bus IDs and angles accordingly.

$wb = $excel.ActiveWorkbook

# --- Transformer sheet ---
$wsT = $wb.Worksheets.Item("Transformer")
$wsT.Range("W23").Value = 0.0001000000013541667
$wsT.Range("X23").Value = 0.00005000000067708333
$wsT.Range("Y23").Value = 0.00005000000067708333
$wsT.Range("W24").Value = 0.0001000000013541667
$wsT.Range("X24").Value = 0.00005000000067708333
$wsT.Range("Y24").Value = 0.00005000000067708333
$wsT.Range("W25").Value = 0.0001000000013541667
$wsT.Range("X25").Value = 0.00005000000067708333
$wsT.Range("Y25").Value = 0.00005000000067708333
$wsT.Range("D26").Value = "xfbus_tpoletop_B"
$wsT.Range("W26").Value = 0.0068
$wsT.Range("X26").Value = 0.006
$wsT.Range("Y26").Value = 0.006
$wsT.Range("D27").Value = "xfbus_tpoletop_B"
$wsT.Range("W27").Value = 0.0068
$wsT.Range("X27").Value = 0.006
$wsT.Range("Y27").Value = 0.006
$wsT.Range("J32").Value = "xfbus_tpoletop_B"
$wsT.Range("W32").Value = 0.01486472334084964
$wsT.Range("X32").Value = 0.01486472334084964
$wsT.Range("Y32").Value = 2.266666666666666
$wsT.Range("Z32").Value = 2.266666666666666

# --- Bus sheet ---
$wsB = $wb.Worksheets.Item("Bus")
$wsB.Range("A3").Value = "632_C"
$wsB.Range("E3").Value = 120
$wsB.Range("A4").Value = "632_A"
$wsB.Range("E4").Value = 0
$wsB.Range("A5").Value = "632_B"
$wsB.Range("E5").Value = -120
$wsB.Range("A6").Value = "633_C"
$wsB.Range("E6").Value = 120
$wsB.Range("A7").Value = "633_A"
$wsB.Range("E7").Value = 0
$wsB.Range("A8").Value = "633_B"
$wsB.Range("E8").Value = -120
$wsB.Range("A9").Value = "634_C"
$wsB.Range("E9").Value = 120
$wsB.Range("A10").Value = "634_A"
$wsB.Range("E10").Value = 0
$wsB.Range("A11").Value = "634_B"
$wsB.Range("E11").Value = -120
$wsB.Range("A12").Value = "645_C"
$wsB.Range("E12").Value = 120
$wsB.Range("A13").Value = "645_B"
$wsB.Range("E13").Value = -120
$wsB.Range("A14").Value = "646_C"
$wsB.Range("E14").Value = 120
$wsB.Range("A15").Value = "646_B"
$wsB.Range("E15").Value = -120
$wsB.Range("A16").Value = "650_C"
$wsB.Range("E16").Value = 120
$wsB.Range("A17").Value = "650_A"
$wsB.Range("E17").Value = 0
$wsB.Range("A18").Value = "650_B"
$wsB.Range("E18").Value = -120
$wsB.Range("A20").Value = "670_C"
$wsB.Range("E20").Value = 120
$wsB.Range("A21").Value = "670_A"
$wsB.Range("E21").Value = 0
$wsB.Range("A22").Value = "670_B"
$wsB.Range("E22").Value = -120
$wsB.Range("A23").Value = "671_C"
$wsB.Range("E23").Value = 120
$wsB.Range("A24").Value = "671_A"
$wsB.Range("E24").Value = 0
$wsB.Range("A25").Value = "671_B"
$wsB.Range("E25").Value = -120
$wsB.Range("A26").Value = "675_C"
$wsB.Range("E26").Value = 120
$wsB.Range("A27").Value = "675_A"
$wsB.Range("E27").Value = 0
$wsB.Range("A28").Value = "675_B"
$wsB.Range("E28").Value = -120
$wsB.Range("A29").Value = "680_C"
$wsB.Range("E29").Value = 120
$wsB.Range("A30").Value = "680_A"
$wsB.Range("E30").Value = 0
$wsB.Range("A31").Value = "680_B"
$wsB.Range("E31").Value = -120
$wsB.Range("A32").Value = "684_C"
$wsB.Range("E32").Value = 120
$wsB.Range("A33").Value = "684_A"
$wsB.Range("E33").Value = 0
$wsB.Range("A34").Value = "692_C"
$wsB.Range("E34").Value = 120
$wsB.Range("A35").Value = "692_A"
$wsB.Range("E35").Value = 0
$wsB.Range("A36").Value = "692_B"
$wsB.Range("E36").Value = -120
$wsB.Range("A37").Value = "brkr_C"
$wsB.Range("E37").Value = 120
$wsB.Range("A38").Value = "brkr_A"
$wsB.Range("E38").Value = 0
$wsB.Range("A39").Value = "brkr_B"
$wsB.Range("E39").Value = -120
$wsB.Range("A42").Value = "mid_C"
$wsB.Range("E42").Value = 120
$wsB.Range("A43").Value = "mid_A"
$wsB.Range("E43").Value = 0
$wsB.Range("A44").Value = "mid_B"
$wsB.Range("E44").Value = -120
$wsB.Range("A45").Value = "rg60_C"
$wsB.Range("E45").Value = 120
$wsB.Range("A46").Value = "rg60_A"
$wsB.Range("E46").Value = 0
$wsB.Range("A47").Value = "rg60_B"
$wsB.Range("E47").Value = -120
$wsB.Range("A48").Value = "sourcebus_C"
$wsB.Range("E48").Value = 120
$wsB.Range("A49").Value = "sourcebus_A"
$wsB.Range("E49").Value = 0
$wsB.Range("A50").Value = "sourcebus_B"
$wsB.Range("E50").Value = -120
$wsB.Range("A52").Value = "xf1_C"
$wsB.Range("E52").Value = 120
$wsB.Range("A53").Value = "xf1_A"
$wsB.Range("E53").Value = 0
$wsB.Range("A54").Value = "xf1_B"
$wsB.Range("E54").Value = -120
$wsB.Range("A55").Value = "xfbus_tpoletop_B"
